# The two leading "helper" rows of the Formulario sheet (an instructions
# row and an example-values row) are removed. Every row below shifts up
# by two: the old header row (row 3) becomes the new row 1 and the old
# sample-data row (row 4) becomes the new row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formulario")

$ws.Range("1:2").EntireRow.Delete()

# The workbook-level named range "Seleción" pointed at Formulario!$A$4:...
# Excel does not auto-adjust it on a plain row delete for this named
# range here, so update it explicitly to track the 2-row shift.
$wb.Names.Item("Seleci" + [char]0x00F3 + "n").RefersTo = "=Formulario!`$A`$2:`$A`$1048576"

# The list validation on B1:B2 referenced INDIRECT(A3); re-point it at
# the new location of the category cell (A1).
$ws.Range("B1:B2").Validation.Delete()
$ws.Range("B1:B2").Validation.Add(3, 1, 1, "INDIRECT(A1)")

# Reflect the final selection state: the whole of the (now two) data rows.
$ws.Range("A1:XFD2").Select()
